$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 159.13333
$ws.Range("I9").Value = 166.27272
$ws.Range("K9").Value = 166.27272
$ws.Range("M9").Value = 2.727280000000007

$ws.Range("H15").Value = 1092.9368
$ws.Range("I15").Value = 1092.9368
$ws.Range("K15").Value = 3278.8104
$ws.Range("M15").Value = -3109.8104

$ws.Range("H51").Value = 10206.667
$ws.Range("I51").Value = 15737.625
$ws.Range("J51").Value = 3885.5715
$ws.Range("K51").Value = 15737.625
$ws.Range("L51").Value = 3885.5715
$ws.Range("M51").Value = -15253.625
$ws.Range("N51").Value = -4853.5715

$ws.Range("H58").Value = 596122.4
$ws.Range("I58").Value = 1050618.9
$ws.Range("J58").Value = 1780.7693
$ws.Range("K58").Value = 3151856.7
$ws.Range("L58").Value = 5342.3079
$ws.Range("M58").Value = -3151706.7
$ws.Range("N58").Value = -5642.3079

$ws.Range("H62").Value = 2328.4285
$ws.Range("I62").Value = 2416.6667
$ws.Range("J62").Value = 1799
$ws.Range("K62").Value = 2416.6667
$ws.Range("L62").Value = 1799
$ws.Range("M62").Value = -1792.6667
$ws.Range("N62").Value = -3047

$ws.Range("H65").Value = 2328.4285
$ws.Range("I65").Value = 2416.6667
$ws.Range("J65").Value = 1799
$ws.Range("K65").Value = 12083.3335
$ws.Range("L65").Value = 8995
$ws.Range("M65").Value = -8963.333500000001
$ws.Range("N65").Value = -15235

$ws.Range("H98").Value = 6945
$ws.Range("I98").Value = 6813.125
$ws.Range("J98").Value = 8000
$ws.Range("K98").Value = 6813.125
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = -5315.125
$ws.Range("N98").Value = -10996

$ws.Range("H113").Value = 112907.11
$ws.Range("I113").Value = 501102.5
$ws.Range("J113").Value = 1994.1428
$ws.Range("K113").Value = 501102.5
$ws.Range("L113").Value = 1994.1428
$ws.Range("M113").Value = -497848.5
$ws.Range("N113").Value = -8502.1428

$ws.Range("H120").Value = 39750.332
$ws.Range("J120").Value = 39750.332
$ws.Range("L120").Value = 39750.332
$ws.Range("N120").Value = -49426.332

$ws.Range("H122").Value = 6945
$ws.Range("I122").Value = 6813.125
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 20439.375
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -17989.375
$ws.Range("N122").Value = -28900

$ws.Range("H125").Value = 1741.091
$ws.Range("I125").Value = 1672.6666
$ws.Range("J125").Value = 1788.4615
$ws.Range("K125").Value = 15053.9994
$ws.Range("L125").Value = 16096.1535
$ws.Range("M125").Value = -12593.9994
$ws.Range("N125").Value = -21016.1535

$ws.Range("H138").Value = 1768.6171
$ws.Range("J138").Value = 3005.5
$ws.Range("L138").Value = 9016.5
$ws.Range("N138").Value = -19296.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1066.33
$ws.Range("I32").Value = 1075.9463
$ws.Range("J32").Value = 938.5714
$ws.Range("K32").Value = 1075.9463
$ws.Range("L32").Value = 938.5714
$ws.Range("M32").Value = -788.9463000000001
$ws.Range("N32").Value = -1512.5714

$ws.Range("H45").Value = 113839.555
$ws.Range("I45").Value = 144942.28
$ws.Range("J45").Value = 4980
$ws.Range("K45").Value = 144942.28
$ws.Range("L45").Value = 4980
$ws.Range("M45").Value = -144565.28
$ws.Range("N45").Value = -5734

$ws.Range("H122").Value = 1425.4286
$ws.Range("I122").Value = 1405.8889
$ws.Range("J122").Value = 1542.6666
$ws.Range("K122").Value = 4217.6667
$ws.Range("L122").Value = 4627.9998
$ws.Range("M122").Value = -1767.6667
$ws.Range("N122").Value = -9527.9998

$ws.Range("H132").Value = 3089.8245
$ws.Range("I132").Value = 3308.796
$ws.Range("J132").Value = 1748.625
$ws.Range("K132").Value = 9926.387999999999
$ws.Range("L132").Value = 5245.875
$ws.Range("M132").Value = -7396.387999999999
$ws.Range("N132").Value = -10305.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 38000
$ws.Range("J2").Value = 38000
$ws.Range("L2").Value = 38000
$ws.Range("N2").Value = -38226

$ws.Range("H13").Value = 35000
$ws.Range("J13").Value = 35000
$ws.Range("L13").Value = 35000
$ws.Range("N13").Value = -35336

$ws.Range("H86").Value = 59010.65
$ws.Range("I86").Value = 77147.07000000001
$ws.Range("J86").Value = 4601.4
$ws.Range("K86").Value = 77147.07000000001
$ws.Range("L86").Value = 4601.4
$ws.Range("M86").Value = -76024.07000000001
$ws.Range("N86").Value = -6847.4

$ws.Range("H89").Value = 59010.65
$ws.Range("I89").Value = 77147.07000000001
$ws.Range("J89").Value = 4601.4
$ws.Range("K89").Value = 385735.35
$ws.Range("L89").Value = 23007
$ws.Range("M89").Value = -380119.35
$ws.Range("N89").Value = -34239

$ws.Range("H105").Value = 62570.91
$ws.Range("I105").Value = 34297.71
$ws.Range("J105").Value = 500805.5
$ws.Range("K105").Value = 34297.71
$ws.Range("L105").Value = 500805.5
$ws.Range("M105").Value = -32550.71
$ws.Range("N105").Value = -504299.5

$ws.Range("H134").Value = 2149.2642
$ws.Range("I134").Value = 1957.1041
$ws.Range("J134").Value = 3994
$ws.Range("K134").Value = 5871.3123
$ws.Range("L134").Value = 11982
$ws.Range("M134").Value = -3336.3123
$ws.Range("N134").Value = -17052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H118").Value = 45741
$ws.Range("J118").Value = 45741
$ws.Range("L118").Value = 45741
$ws.Range("N118").Value = -49055

$ws.Range("H132").Value = 3563.0527
$ws.Range("I132").Value = 3309.8667
$ws.Range("J132").Value = 4512.5
$ws.Range("K132").Value = 9929.6001
$ws.Range("L132").Value = 13537.5
$ws.Range("M132").Value = -7399.6001
$ws.Range("N132").Value = -18597.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 530.6896400000001
$ws.Range("I113").Value = 536
$ws.Range("J113").Value = 527.44446
$ws.Range("K113").Value = 1608
$ws.Range("L113").Value = 1582.33338
$ws.Range("M113").Value = 562
$ws.Range("N113").Value = -5922.33338

$ws.Range("H131").Value = 1289.5657
$ws.Range("I131").Value = 915.5714
$ws.Range("J131").Value = 1318.0217
$ws.Range("K131").Value = 2746.7142
$ws.Range("L131").Value = 3954.0651
$ws.Range("M131").Value = 2293.2858
$ws.Range("N131").Value = -14034.0651

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2502250
$ws.Range("I5").Value = 5000000
$ws.Range("K5").Value = 5000000
$ws.Range("M5").Value = -4999888

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H126").Value = 4769.5713
$ws.Range("I126").Value = 4931.1665
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 14793.4995
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -12323.4995
$ws.Range("N126").Value = -16340

$ws.Range("H132").Value = 1858.597
$ws.Range("I132").Value = 1898.1864
$ws.Range("J132").Value = 1566.625
$ws.Range("K132").Value = 5694.5592
$ws.Range("L132").Value = 4699.875
$ws.Range("M132").Value = -3164.5592
$ws.Range("N132").Value = -9759.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 57910
$ws.Range("I40").Value = 144982.86
$ws.Range("K40").Value = 144982.86
$ws.Range("M40").Value = -144846.86

$ws.Range("H60").Value = 49990
$ws.Range("J60").Value = 49990
$ws.Range("L60").Value = 49990
$ws.Range("N60").Value = -51008

$ws.Range("H122").Value = 2584.923
$ws.Range("I122").Value = 2509.4546
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7528.3638
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5078.3638
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 2752.027
$ws.Range("I132").Value = 2846.0908
$ws.Range("K132").Value = 8538.2724
$ws.Range("M132").Value = -6008.2724

$ws.Range("H136").Value = 1244.738
$ws.Range("I136").Value = 1190.8055
$ws.Range("J136").Value = 1568.3334
$ws.Range("K136").Value = 3572.4165
$ws.Range("L136").Value = 4705.0002
$ws.Range("M136").Value = -1022.4165
$ws.Range("N136").Value = -9805.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3042.4
$ws.Range("I126").Value = 3461.3333
$ws.Range("J126").Value = 2414
$ws.Range("K126").Value = 10383.9999
$ws.Range("L126").Value = 7242
$ws.Range("M126").Value = -7913.999899999999
$ws.Range("N126").Value = -12182

$ws.Range("H132").Value = 1655.8667
$ws.Range("I132").Value = 1752.4067
$ws.Range("J132").Value = 1299.875
$ws.Range("K132").Value = 5257.2201
$ws.Range("L132").Value = 3899.625
$ws.Range("M132").Value = -2727.2201
$ws.Range("N132").Value = -8959.625

$ws.Range("H136").Value = 581.69385
$ws.Range("J136").Value = 1625.75
$ws.Range("L136").Value = 4877.25
$ws.Range("N136").Value = -9977.25
